{"js": "// Apply the Nmap/Metasploit log timestamp + reordering edits described\n// by the diff. We operate on Word.Paragraph.text (Office.js represents\n// <w:br/> as the vertical-tab char \"\\u000b\") so we can do precise,\n// localized substring replacements without disturbing any other content\n// in these long single-run \"console output\" paragraphs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nfunction applyReplacements(paragraph, replacements) {\n  // replacements: array of [oldSubstring, newSubstring]\n  let text = paragraph.text;\n  let changed = false;\n  for (const [oldSub, newSub] of replacements) {\n    if (text.indexOf(oldSub) === -1) {\n      throw new Error(\"Substring not found: \" + oldSub);\n    }\n    text = text.split(oldSub).join(newSub);\n    changed = true;\n  }\n  if (changed) {\n    paragraph.insertText(text, Word.InsertLocation.replace);\n  }\n}\n\nconst BR = \"\\u000b\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const t = p.text;\n  if (!t) continue;\n\n  // 1) Nmap scan block: initiated time, first host latency, and the\n  //    \"Nmap done\" summary line (time + scan duration).\n  if (t.indexOf(\"# Nmap 7.80 scan initiated\") !== -1) {\n    applyReplacements(p, [\n      [\n        \"# Nmap 7.80 scan initiated Thu Jul 11 10:32:56 2024 as:\",\n        \"# Nmap 7.80 scan initiated Thu Jul 11 10:36:53 2024 as:\",\n      ],\n      [\n        \"Host is up (0.00014s latency).\",\n        \"Host is up (0.00012s latency).\",\n      ],\n      [\n        \"# Nmap done at Thu Jul 11 10:33:04 2024 -- 2 IP addresses (2 hosts up) scanned in 8.68 seconds\",\n        \"# Nmap done at Thu Jul 11 10:37:01 2024 -- 2 IP addresses (2 hosts up) scanned in 7.99 seconds\",\n      ],\n    ]);\n    continue;\n  }\n\n  // 2) Metasploit exploitation log: move the \"Sending stage\" /\n  //    \"Meterpreter session opened\" lines so they appear right after\n  //    \"Enumerating local_data_id\" (before \"Found exploitable\" /\n  //    \"Command Stager progress\"), and bump the session port + timestamp.\n  //    Also update the later session-table line that repeats the port.\n  if (t.indexOf(\"[*] Enumerating local_data_id values for host_id 1\") !== -1) {\n    const oldBlock =\n      \"[+] Found exploitable local_data_id 15 for host_id 1\" + BR +\n      \"[*] Command Stager progress - 100.00% done (1118/1118 bytes)\" + BR +\n      \"[*] Sending stage (1017704 bytes) to 10.33.102.225\" + BR +\n      \"[*] Meterpreter session 1 opened (10.33.102.224:4444 -> 10.33.102.225:35958) at 2024-07-11 10:34:19 +0700\";\n    const newBlock =\n      \"[*] Sending stage (1017704 bytes) to 10.33.102.225\" + BR +\n      \"[*] Meterpreter session 1 opened (10.33.102.224:4444 -> 10.33.102.225:40740) at 2024-07-11 10:38:02 +0700\" + BR +\n      \"[+] Found exploitable local_data_id 15 for host_id 1\" + BR +\n      \"[*] Command Stager progress - 100.00% done (1118/1118 bytes)\";\n\n    applyReplacements(p, [\n      [oldBlock, newBlock],\n      [\n        \"            x                                             10.33.102.225:35958\",\n        \"            x                                             10.33.102.225:40740\",\n      ],\n    ]);\n    continue;\n  }\n}\n", "ps1": "# Apply the Nmap/Metasploit log timestamp + reordering edits described\n# by the diff, using Word COM interop against $word.ActiveDocument.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $found = $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n    if (-not $found) {\n        throw \"Text not found: $find\"\n    }\n}\n\n# 1) Nmap scan block: initiated time, first host latency, and the\n#    \"Nmap done\" summary line (time + scan duration).\nReplace-Text \"# Nmap 7.80 scan initiated Thu Jul 11 10:32:56 2024 as:\" \"# Nmap 7.80 scan initiated Thu Jul 11 10:36:53 2024 as:\"\nReplace-Text \"Host is up (0.00014s latency).\" \"Host is up (0.00012s latency).\"\nReplace-Text \"# Nmap done at Thu Jul 11 10:33:04 2024 -- 2 IP addresses (2 hosts up) scanned in 8.68 seconds\" \"# Nmap done at Thu Jul 11 10:37:01 2024 -- 2 IP addresses (2 hosts up) scanned in 7.99 seconds\"\n\n# 2) Metasploit exploitation log: move the \"Sending stage\" /\n#    \"Meterpreter session opened\" lines so they appear right after\n#    \"Enumerating local_data_id\" (before \"Found exploitable\" /\n#    \"Command Stager progress\"), and bump the session port + timestamp.\n#    Also update the later session-table line that repeats the port.\n$BR = [char]11\n\nforeach ($p in $d.Paragraphs) {\n    $rngFull = $p.Range\n    $full = $rngFull.Text\n    if ($full -and $full.Contains(\"[*] Enumerating local_data_id values for host_id 1\")) {\n        # Drop the trailing paragraph-mark char so assigning .Text back\n        # doesn't split off a new empty paragraph.\n        $t = $full.Substring(0, $full.Length - 1)\n\n        $oldBlock = \"[+] Found exploitable local_data_id 15 for host_id 1\" + $BR + `\n                    \"[*] Command Stager progress - 100.00% done (1118/1118 bytes)\" + $BR + `\n                    \"[*] Sending stage (1017704 bytes) to 10.33.102.225\" + $BR + `\n                    \"[*] Meterpreter session 1 opened (10.33.102.224:4444 -> 10.33.102.225:35958) at 2024-07-11 10:34:19 +0700\"\n        $newBlock = \"[*] Sending stage (1017704 bytes) to 10.33.102.225\" + $BR + `\n                    \"[*] Meterpreter session 1 opened (10.33.102.224:4444 -> 10.33.102.225:40740) at 2024-07-11 10:38:02 +0700\" + $BR + `\n                    \"[+] Found exploitable local_data_id 15 for host_id 1\" + $BR + `\n                    \"[*] Command Stager progress - 100.00% done (1118/1118 bytes)\"\n\n        if (-not $t.Contains($oldBlock)) {\n            throw \"Expected meterpreter log block not found\"\n        }\n        $newText = $t.Replace($oldBlock, $newBlock)\n        $newText = $newText.Replace(\n            \"            x                                             10.33.102.225:35958\",\n            \"            x                                             10.33.102.225:40740\")\n\n        $rng = $d.Range($rngFull.Start, $rngFull.End - 1)\n        $rng.Text = $newText\n        break\n    }\n}\n"}
